$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing the existing rows 25-49 down to 26-50.
$ws.Range("A25").EntireRow.Insert()

# Populate the newly inserted row 25 with the latest week's reading
# (matches the style of the surrounding "Pepino dulce" / Macroferia Regional de Talca rows).
$ws.Range("A25").Value2 = 5
$ws.Range("B25").Value2 = "Macroferia Regional de Talca"
$ws.Range("C25").Value2 = "Maule"
$ws.Range("D25").Value2 = 44804
$ws.Range("E25").Value2 = 7
$ws.Range("F25").Value2 = 100112043
$ws.Range("G25").Value2 = "Pepino dulce"
$ws.Range("H25").Value2 = "Cultivar IV Región"
$ws.Range("I25").Value2 = "Primera"
$ws.Range("J25").Value2 = 300
$ws.Range("K25").Value2 = 15000
$ws.Range("L25").Value2 = 15000
$ws.Range("M25").Value2 = 15000
$ws.Range("N25").Value2 = "$/bandeja 18 kilos"
$ws.Range("O25").Value2 = "Provincia de Limarí"
$ws.Range("P25").Value2 = 833
$ws.Range("Q25").Value2 = 18
$ws.Range("R25").Value2 = "Hortaliza"
